# Update time_taken (F column) timestamps on the "data" sheet (rows 2-82)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Cells.Item(2, 6).Value = "2021-10-05 14:33:47.458487"
$ws.Cells.Item(3, 6).Value = "2021-10-05 14:33:47.458495"
$ws.Cells.Item(4, 6).Value = "2021-10-05 14:33:47.458498"
$ws.Cells.Item(5, 6).Value = "2021-10-05 14:33:47.458501"
$ws.Cells.Item(6, 6).Value = "2021-10-05 14:33:47.458504"
$ws.Cells.Item(7, 6).Value = "2021-10-05 14:33:47.458506"
$ws.Cells.Item(8, 6).Value = "2021-10-05 14:33:47.458509"
$ws.Cells.Item(9, 6).Value = "2021-10-05 14:33:47.458511"
$ws.Cells.Item(10, 6).Value = "2021-10-05 14:33:47.458514"
$ws.Cells.Item(11, 6).Value = "2021-10-05 14:33:47.458517"
$ws.Cells.Item(12, 6).Value = "2021-10-05 14:33:47.458519"
$ws.Cells.Item(13, 6).Value = "2021-10-05 14:33:47.458522"
$ws.Cells.Item(14, 6).Value = "2021-10-05 14:33:47.458524"
$ws.Cells.Item(15, 6).Value = "2021-10-05 14:33:47.458527"
$ws.Cells.Item(16, 6).Value = "2021-10-05 14:33:47.458529"
$ws.Cells.Item(17, 6).Value = "2021-10-05 14:33:47.458531"
$ws.Cells.Item(18, 6).Value = "2021-10-05 14:33:47.458534"
$ws.Cells.Item(19, 6).Value = "2021-10-05 14:33:47.458539"
$ws.Cells.Item(20, 6).Value = "2021-10-05 14:33:47.458542"
$ws.Cells.Item(21, 6).Value = "2021-10-05 14:33:47.458544"
$ws.Cells.Item(22, 6).Value = "2021-10-05 14:33:47.458546"
$ws.Cells.Item(23, 6).Value = "2021-10-05 14:33:47.458549"
$ws.Cells.Item(24, 6).Value = "2021-10-05 14:33:47.458553"
$ws.Cells.Item(25, 6).Value = "2021-10-05 14:33:47.458555"
$ws.Cells.Item(26, 6).Value = "2021-10-05 14:33:47.458558"
$ws.Cells.Item(27, 6).Value = "2021-10-05 14:33:47.458560"
$ws.Cells.Item(28, 6).Value = "2021-10-05 14:33:47.458563"
$ws.Cells.Item(29, 6).Value = "2021-10-05 14:33:47.458565"
$ws.Cells.Item(30, 6).Value = "2021-10-05 14:33:47.458568"
$ws.Cells.Item(31, 6).Value = "2021-10-05 14:33:47.458570"
$ws.Cells.Item(32, 6).Value = "2021-10-05 14:33:47.458573"
$ws.Cells.Item(33, 6).Value = "2021-10-05 14:33:47.458575"
$ws.Cells.Item(34, 6).Value = "2021-10-05 14:33:47.458578"
$ws.Cells.Item(35, 6).Value = "2021-10-05 14:33:47.458581"
$ws.Cells.Item(36, 6).Value = "2021-10-05 14:33:47.458583"
$ws.Cells.Item(37, 6).Value = "2021-10-05 14:33:47.458585"
$ws.Cells.Item(38, 6).Value = "2021-10-05 14:33:47.458588"
$ws.Cells.Item(39, 6).Value = "2021-10-05 14:33:47.458590"
$ws.Cells.Item(40, 6).Value = "2021-10-05 14:33:47.458593"
$ws.Cells.Item(41, 6).Value = "2021-10-05 14:33:47.458595"
$ws.Cells.Item(42, 6).Value = "2021-10-05 14:33:47.458598"
$ws.Cells.Item(43, 6).Value = "2021-10-05 14:33:47.458601"
$ws.Cells.Item(44, 6).Value = "2021-10-05 14:33:47.458603"
$ws.Cells.Item(45, 6).Value = "2021-10-05 14:33:47.458606"
$ws.Cells.Item(46, 6).Value = "2021-10-05 14:33:47.458608"
$ws.Cells.Item(47, 6).Value = "2021-10-05 14:33:47.458610"
$ws.Cells.Item(48, 6).Value = "2021-10-05 14:33:47.458613"
$ws.Cells.Item(49, 6).Value = "2021-10-05 14:33:47.458615"
$ws.Cells.Item(50, 6).Value = "2021-10-05 14:33:47.458618"
$ws.Cells.Item(51, 6).Value = "2021-10-05 14:33:47.458620"
$ws.Cells.Item(52, 6).Value = "2021-10-05 14:33:47.458623"
$ws.Cells.Item(53, 6).Value = "2021-10-05 14:33:47.458625"
$ws.Cells.Item(54, 6).Value = "2021-10-05 14:33:47.458628"
$ws.Cells.Item(55, 6).Value = "2021-10-05 14:33:47.458630"
$ws.Cells.Item(56, 6).Value = "2021-10-05 14:33:47.458633"
$ws.Cells.Item(57, 6).Value = "2021-10-05 14:33:47.458635"
$ws.Cells.Item(58, 6).Value = "2021-10-05 14:33:47.458638"
$ws.Cells.Item(59, 6).Value = "2021-10-05 14:33:47.458640"
$ws.Cells.Item(60, 6).Value = "2021-10-05 14:33:47.458643"
$ws.Cells.Item(61, 6).Value = "2021-10-05 14:33:47.458645"
$ws.Cells.Item(62, 6).Value = "2021-10-05 14:33:47.458647"
$ws.Cells.Item(63, 6).Value = "2021-10-05 14:33:47.458650"
$ws.Cells.Item(64, 6).Value = "2021-10-05 14:33:47.458652"
$ws.Cells.Item(65, 6).Value = "2021-10-05 14:33:47.458655"
$ws.Cells.Item(66, 6).Value = "2021-10-05 14:33:47.458658"
$ws.Cells.Item(67, 6).Value = "2021-10-05 14:33:47.458661"
$ws.Cells.Item(68, 6).Value = "2021-10-05 14:33:47.458663"
$ws.Cells.Item(69, 6).Value = "2021-10-05 14:33:47.458666"
$ws.Cells.Item(70, 6).Value = "2021-10-05 14:33:47.458668"
$ws.Cells.Item(71, 6).Value = "2021-10-05 14:33:47.458670"
$ws.Cells.Item(72, 6).Value = "2021-10-05 14:33:47.458673"
$ws.Cells.Item(73, 6).Value = "2021-10-05 14:33:47.458675"
$ws.Cells.Item(74, 6).Value = "2021-10-05 14:33:47.458678"
$ws.Cells.Item(75, 6).Value = "2021-10-05 14:33:47.458680"
$ws.Cells.Item(76, 6).Value = "2021-10-05 14:33:47.458683"
$ws.Cells.Item(77, 6).Value = "2021-10-05 14:33:47.458685"
$ws.Cells.Item(78, 6).Value = "2021-10-05 14:33:47.458690"
$ws.Cells.Item(79, 6).Value = "2021-10-05 14:33:47.458692"
$ws.Cells.Item(80, 6).Value = "2021-10-05 14:33:47.458695"
$ws.Cells.Item(81, 6).Value = "2021-10-05 14:33:47.458697"
$ws.Cells.Item(82, 6).Value = "2021-10-05 14:33:47.458700"


# Add a new "metadata" sheet right after "data", mirroring the panelapp
# scraper's per-panel metadata record (data_name/data_id/data_version/...).
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Reuse the bold/centered/bordered header style already used on "data" (style id 1)
# by copy/pasting formats from an existing styled cell, instead of re-building
# the format from scratch (which would create a near-duplicate style entry).
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Early-onset Dementia"
$meta.Range("C2").Value = 24
# Force "0.148" to be stored as text (matches source data), not coerced to a
# number: write it as a text-producing formula, then copy/paste-values over
# itself to flatten to a static value without picking up a quote-prefix
# style (keeps the cell styleless, like the target file).
$meta.Range("D2").Formula = '="0.148"'
$meta.Range("D2").Copy()
$meta.Range("D2").PasteSpecial(-4163)
$meta.Range("E2").Value = "2021-08-31T02:46:36.806048Z"
$meta.Range("F2").Value = "2021-10-05 14:33:47.454899"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/24/?format=json"

# Keep "data" as the active/selected sheet (unchanged in the target workbook).
$ws.Activate()
